$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("projects")
$ws = $wb.Worksheets.Item("terms")

# --- Update selection on the "projects" sheet (it stays the inactive tab) ---
$ws1.Range("F1").Select() | Out-Null

$ws.Activate() | Out-Null

# --- Fix the shared string text used for the "terms" header (budgetedAmount -> Budgeted Amount) ---
$ws.Range("B1").Value = "Budgeted Amount"

# --- Insert new "Item Work Type Reference" column (becomes column B) ---
$ws.Columns("B:B").Insert()
$ws.Range("B1").Value = "Item Work Type Reference"
$ws.Columns("B").ColumnWidth = 24.830729166666668
$ws.Columns("C").ColumnWidth = 18.498697916666668

$ws.Range("B2").Value = "OTHER"
$ws.Range("B3").Value = "OTHER"
$ws.Range("B4").Value = "WORKS"

# --- Split the old row 4 (20000 / 43282-43373) into a WORKS/5 row and a WORKS/15000 row ---
$ws.Rows("5:5").Insert()

$ws.Range("A5").Value = "GB01"
$ws.Range("B5").Value = "WORKS"
$ws.Range("C5").Value = 15000
$ws.Range("D5").Value = 43282
$ws.Range("E5").Value = 43373

$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 43101
$ws.Range("E4").Value = 43190

# --- Selection to match the saved view state (terms stays the active/selected tab) ---
$ws.Range("C3").Select() | Out-Null
